# JSON input card functionality
# The author's date-stamped "last edited" footer field moved on by a day:
# every datetimeFigureOut placeholder (slide master + all slide layouts)
# gets its cached text bumped from 25/01/2021 -> 26/01/2021.

$p = $ppt.ActivePresentation
$newDate = "26/01/2021"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master's own Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout hanging off the master.
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Update-DatePlaceholder $layout.Shapes
}
